{"js": "// The diff inserts four new empty paragraphs (style \"Normal\", justified\n// \"both\", single empty run) immediately after the paragraph ending in\n// \"...ayudan a encontrar los problemas lo antes posible.\" and before the\n// pre-existing blank paragraph that precedes the \"Anexos\" heading.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its distinctive trailing text.\nconst marker = \"ayudan a encontrar los problemas lo antes posible.\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text && p.text.indexOf(marker) !== -1) {\n    anchor = p;\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Anchor paragraph not found\");\n}\n\n// Insert four empty paragraphs right after the anchor paragraph. Each new\n// paragraph inherits the anchor's paragraph formatting (style \"Normal\",\n// justified \"both\"), matching the diff. Chaining off the newly created\n// paragraph each time keeps them directly following the anchor and before\n// the following (pre-existing) blank paragraph.\nlet current = anchor;\nfor (let i = 0; i < 4; i++) {\n  current = current.insertParagraph(\"\", \"After\");\n}\n\nawait context.sync();\n", "ps1": "# The diff inserts four new empty paragraphs (style \"Normal\", justified\n# \"both\", single empty run) immediately after the paragraph ending in\n# \"...ayudan a encontrar los problemas lo antes posible.\" and before the\n# pre-existing blank paragraph that precedes the \"Anexos\" heading.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph by its distinctive trailing text.\n$marker = \"ayudan a encontrar los problemas lo antes posible.\"\n$anchorRange = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$marker*\") {\n        $anchorRange = $p.Range\n        break\n    }\n}\n\nif ($anchorRange -eq $null) {\n    throw \"Anchor paragraph not found\"\n}\n\n# Insert four empty paragraphs right after the anchor paragraph. Each call\n# to InsertParagraphAfter() on the (unmoved) anchor range inherits the\n# anchor's paragraph formatting (style \"Normal\", justified \"both\") and\n# inserts immediately after the anchor paragraph -- ahead of any paragraph\n# already inserted by a previous iteration -- so the four new paragraphs\n# end up directly after the anchor, in order, and before the following\n# (pre-existing) blank paragraph, exactly matching the diff.\nfor ($i = 0; $i -lt 4; $i++) {\n    $anchorRange.InsertParagraphAfter()\n}\n"}
